$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("臨床イベント")

# Update rows 2-7 with new clinical event data
$data = @(
    @("2024-06-20", "低血糖", "血糖52mg/dL、冷汗・手指振戦、補食で回復"),
    @("2024-07-05", "低血糖", "血糖48mg/dL、運動後、ブドウ糖摂取"),
    @("2024-07-25", "高血糖", "血糖320mg/dL、インスリン打ち忘れ"),
    @("2024-08-10", "低血糖", "血糖55mg/dL、夕食遅延時"),
    @("2024-09-15", "低血糖", "血糖60mg/dL、体育後、軽度"),
    @("2024-11-20", "高血糖", "血糖280mg/dL、感冒時sick day")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    # Force the date-looking string to stay text (not auto-converted to a date
    # serial), then clear the format again so no extra cell style is introduced.
    $ws.Cells.Item($row, 1).NumberFormat = "@"
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 1).ClearFormats()
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}

# Remove now-unused rows 8 and 9 (dimension shrinks from A1:C9 to A1:C7)
$ws.Range("A8:C9").Delete()
